$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right after the header row (row 1), pushing all
# existing data (old rows 2..205) down to rows 6..209.
$ws.Rows("2:5").Insert()

# Fill in the newly inserted rows with the Twilight-themed songs.
$ws.Range("B2").Value = "Edward Cullen - Bella's Lullaby"
$ws.Range("C2").Value = "https://www.youtube.com/watch?v=zQME-ChSwNM"

$ws.Range("B3").Value = "Alexandre Desplat - New Moon (The Meadow)"
$ws.Range("C3").Value = "https://www.youtube.com/watch?v=7kY5bQEU5gQ"

$ws.Range("B4").Value = "New Moon - Edward leaves - Alexandre Desplat"
$ws.Range("C4").Value = "https://www.youtube.com/watch?v=MS4Tf9mr44M"

$ws.Range("B5").Value = "New Moon OST - Dreamcatcher - Alexandre Desplat"
$ws.Range("C5").Value = "https://www.youtube.com/watch?v=39Kvcgug2J0"
